$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 1056.4286
$ws.Range("I46").Value = 400
$ws.Range("J46").Value = 1165.8334
$ws.Range("K46").Value = 1200
$ws.Range("L46").Value = 3497.5002
$ws.Range("M46").Value = -1081
$ws.Range("N46").Value = -3735.5002
$ws.Range("H59").Value = 1038.6364
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 1038.6364
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 3115.9092
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -4229.9092
$ws.Range("H60").Value = 1056.4286
$ws.Range("I60").Value = 400
$ws.Range("J60").Value = 1165.8334
$ws.Range("K60").Value = 1200
$ws.Range("L60").Value = 3497.5002
$ws.Range("M60").Value = -716
$ws.Range("N60").Value = -4465.5002
$ws.Range("H138").Value = 2631.84
$ws.Range("I138").Value = 1470.3043
$ws.Range("J138").Value = 2978.7922
$ws.Range("K138").Value = 4410.9129
$ws.Range("L138").Value = 8936.3766
$ws.Range("M138").Value = 729.0870999999997
$ws.Range("N138").Value = -19216.3766

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1199.8
$ws.Range("I4").Value = 750
$ws.Range("J4").Value = 1499.6666
$ws.Range("K4").Value = 750
$ws.Range("L4").Value = 1499.6666
$ws.Range("M4").Value = -634
$ws.Range("N4").Value = -1731.6666
$ws.Range("H32").Value = 10068.976
$ws.Range("I32").Value = 6501.1606
$ws.Range("J32").Value = 17468.889
$ws.Range("K32").Value = 6501.1606
$ws.Range("L32").Value = 17468.889
$ws.Range("M32").Value = -6214.1606
$ws.Range("N32").Value = -18042.889
$ws.Range("H35").Value = 24533
$ws.Range("J35").Value = 34221.668
$ws.Range("L35").Value = 34221.668
$ws.Range("N35").Value = -35033.668
$ws.Range("H80").Value = 36124.777
$ws.Range("J80").Value = 36124.777
$ws.Range("L80").Value = 36124.777
$ws.Range("N80").Value = -38120.777
$ws.Range("H83").Value = 36124.777
$ws.Range("J83").Value = 36124.777
$ws.Range("L83").Value = 108374.331
$ws.Range("N83").Value = -118358.331
$ws.Range("H97").Value = 1198.36
$ws.Range("I97").Value = 945.45
$ws.Range("J97").Value = 2210
$ws.Range("K97").Value = 945.45
$ws.Range("L97").Value = 2210
$ws.Range("M97").Value = -449.45
$ws.Range("N97").Value = -3202

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 2879682.8
$ws.Range("I7").Value = 6683334.5
$ws.Range("J7").Value = 26944.25
$ws.Range("K7").Value = 6683334.5
$ws.Range("L7").Value = 26944.25
$ws.Range("M7").Value = -6683221.5
$ws.Range("N7").Value = -27170.25
$ws.Range("H20").Value = 9528.444
$ws.Range("I20").Value = 1729.2727
$ws.Range("J20").Value = 21784.285
$ws.Range("K20").Value = 1729.2727
$ws.Range("L20").Value = 21784.285
$ws.Range("M20").Value = -1482.2727
$ws.Range("N20").Value = -22278.285
$ws.Range("H82").Value = 15419.214
$ws.Range("I82").Value = 2596
$ws.Range("J82").Value = 32516.834
$ws.Range("K82").Value = 2596
$ws.Range("L82").Value = 32516.834
$ws.Range("M82").Value = -2213
$ws.Range("N82").Value = -33282.834
$ws.Range("H85").Value = 15419.214
$ws.Range("I85").Value = 2596
$ws.Range("J85").Value = 32516.834
$ws.Range("K85").Value = 2596
$ws.Range("L85").Value = 32516.834
$ws.Range("M85").Value = -1270
$ws.Range("N85").Value = -35168.834
$ws.Range("H99").Value = 3636.5
$ws.Range("I99").Value = 1320
$ws.Range("K99").Value = 1320
$ws.Range("M99").Value = 178

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 3753
$ws.Range("H31").Value = 2971.4822
$ws.Range("I31").Value = 1338.3235
$ws.Range("J31").Value = 5495.4546
$ws.Range("K31").Value = 1338.3235
$ws.Range("L31").Value = 5495.4546
$ws.Range("M31").Value = -1043.3235
$ws.Range("N31").Value = -6085.4546
$ws.Range("H34").Value = 2971.4822
$ws.Range("I34").Value = 1338.3235
$ws.Range("J34").Value = 5495.4546
$ws.Range("K34").Value = 1338.3235
$ws.Range("L34").Value = 5495.4546
$ws.Range("M34").Value = -1136.3235
$ws.Range("N34").Value = -5899.4546
$ws.Range("H109").Value = 29999.5
$ws.Range("J109").Value = 29999.5
$ws.Range("L109").Value = 29999.5
$ws.Range("N109").Value = -32079.5
$ws.Range("H122").Value = 2600.682
$ws.Range("I122").Value = 1890.3125
$ws.Range("J122").Value = 4495
$ws.Range("K122").Value = 5670.9375
$ws.Range("L122").Value = 13485
$ws.Range("M122").Value = -3220.9375
$ws.Range("N122").Value = -18385

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3449279.2
$ws.Range("J4").Value = 8845
$ws.Range("L4").Value = 26535
$ws.Range("N4").Value = -26759
$ws.Range("H131").Value = 14286670
$ws.Range("J131").Value = 1147.6
$ws.Range("L131").Value = 3442.8
$ws.Range("N131").Value = -13522.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 29000
$ws.Range("J4").Value = 29000
$ws.Range("L4").Value = 29000
$ws.Range("N4").Value = -29224
$ws.Range("H70").Value = 6518.2
$ws.Range("I70").Value = 5825.278
$ws.Range("J70").Value = 8300
$ws.Range("K70").Value = 5825.278
$ws.Range("L70").Value = 8300
$ws.Range("M70").Value = -5555.278
$ws.Range("N70").Value = -8840
$ws.Range("H73").Value = 6518.2
$ws.Range("I73").Value = 5825.278
$ws.Range("J73").Value = 8300
$ws.Range("K73").Value = 5825.278
$ws.Range("L73").Value = 8300
$ws.Range("M73").Value = -4889.278
$ws.Range("N73").Value = -10172

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4486.6665
$ws.Range("I7").Value = 3300
$ws.Range("J7").Value = 6266.6665
$ws.Range("K7").Value = 3300
$ws.Range("L7").Value = 6266.6665
$ws.Range("M7").Value = -3188
$ws.Range("N7").Value = -6490.6665
$ws.Range("H22").Value = 8335091
$ws.Range("I22").Value = 14707186
$ws.Range("J22").Value = 2350.923
$ws.Range("K22").Value = 14707186
$ws.Range("L22").Value = 2350.923
$ws.Range("M22").Value = -14706891
$ws.Range("N22").Value = -2940.923
$ws.Range("H27").Value = 8335091
$ws.Range("I27").Value = 14707186
$ws.Range("J27").Value = 2350.923
$ws.Range("K27").Value = 14707186
$ws.Range("L27").Value = 2350.923
$ws.Range("M27").Value = -14707079
$ws.Range("N27").Value = -2564.923
$ws.Range("H61").Value = 2332.4666
$ws.Range("I61").Value = 2165.5833
$ws.Range("K61").Value = 2165.5833
$ws.Range("M61").Value = -1963.5833
$ws.Range("H68").Value = 868.561
$ws.Range("I68").Value = 734.12823
$ws.Range("K68").Value = 734.12823
$ws.Range("M68").Value = 14.87176999999997
$ws.Range("H71").Value = 868.561
$ws.Range("I71").Value = 734.12823
$ws.Range("K71").Value = 3670.64115
$ws.Range("M71").Value = 73.35884999999962
$ws.Range("H75").Value = 44500
$ws.Range("J75").Value = 44500
$ws.Range("L75").Value = 44500
$ws.Range("N75").Value = -46372
$ws.Range("H78").Value = 44500
$ws.Range("J78").Value = 44500
$ws.Range("L78").Value = 133500
$ws.Range("N78").Value = -142860
$ws.Range("H113").Value = 2332.4666
$ws.Range("I113").Value = 2165.5833
$ws.Range("K113").Value = 2165.5833
$ws.Range("M113").Value = 4.416700000000219
$ws.Range("H126").Value = 4486.6665
$ws.Range("I126").Value = 3300
$ws.Range("J126").Value = 6266.6665
$ws.Range("K126").Value = 9900
$ws.Range("L126").Value = 18799.9995
$ws.Range("M126").Value = -7430
$ws.Range("N126").Value = -23739.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 16166.667
$ws.Range("J86").Value = 16166.667
$ws.Range("L86").Value = 16166.667
$ws.Range("N86").Value = -18412.667
$ws.Range("H89").Value = 16166.667
$ws.Range("J89").Value = 16166.667
$ws.Range("L89").Value = 80833.33499999999
$ws.Range("N89").Value = -92065.33499999999
$ws.Range("H94").Value = 34500
$ws.Range("J94").Value = 34500
$ws.Range("L94").Value = 34500
$ws.Range("N94").Value = -36302
